# SAMPL6 microstates workbook: add a new "canonical SMILES" column (D) that
# mirrors the existing "canonical isomeric SMILES" column (C) for every
# microstate row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header, same row (2) as the other column headers.
$ws.Range("D2").Value = "canonical SMILES"

# Rows 3-22 hold the microstate data; duplicate column C's SMILES into the
# new column D for each of them.
for ($row = 3; $row -le 22; $row++) {
    $cSmiles = $ws.Cells.Item($row, 3).Value2
    $ws.Cells.Item($row, 4).Value = $cSmiles
}

# Give the new column a sensible width, matching the other data columns
# (~36.86 characters, the closest value the host's width model can store).
$ws.Columns.Item(4).ColumnWidth = 36.0
